$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "260.99"
Set-TextValue "E2" "0.93%"
Set-TextValue "D3" "27.05"
Set-TextValue "E3" "0.79%"
Set-TextValue "E4" "1.30%"
Set-TextValue "D5" "0.06175"
Set-TextValue "E5" "3.56%"
Set-TextValue "D6" "6.684"
Set-TextValue "E6" "0.60%"
Set-TextValue "E7" "-0.81%"
Set-TextValue "D8" "0.9145"
Set-TextValue "E8" "-0.83%"
Set-TextValue "E9" "1.41%"
Set-TextValue "D10" "0.04675"
Set-TextValue "E10" "8.07%"
Set-TextValue "D11" "0.07088"
Set-TextValue "E11" "1.15%"
Set-TextValue "D12" "0.03150"
Set-TextValue "E12" "3.15%"
Set-TextValue "D13" "0.09039"
Set-TextValue "E13" "-0.81%"
Set-TextValue "D14" "0.001535"
Set-TextValue "E14" "0.72%"
Set-TextValue "D15" "0.0006147"
Set-TextValue "E15" "1.52%"
Set-TextValue "D16" "0.006020"
Set-TextValue "E16" "-0.95%"
Set-TextValue "D17" "3.450"
Set-TextValue "E17" "0.00%"
Set-TextValue "D18" "3.166"
Set-TextValue "E18" "0.53%"
Set-TextValue "D19" "2.158"
Set-TextValue "E19" "0.19%"
Set-TextValue "E20" "-0.88%"
Set-TextValue "D21" "0.1300"
Set-TextValue "E21" "0.95%"
Set-TextValue "D22" "4.096"
Set-TextValue "E22" "2.29%"
Set-TextValue "D23" "0.04222"
Set-TextValue "E23" "-0.18%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "-0.04%"
Set-TextValue "E25" "-5.67%"
Set-TextValue "E26" "0.12%"
Set-TextValue "D40" "0.03877"
Set-TextValue "E40" "1.32%"
Set-TextValue "E41" "-0.11%"
Set-TextValue "D42" "0.004092"
Set-TextValue "E42" "8.95%"
Set-TextValue "E43" "8.08%"
Set-TextValue "E44" "-10.05%"
Set-TextValue "D45" "0.00005156"
Set-TextValue "E45" "0.62%"
Set-TextValue "E46" "0.12%"
Set-TextValue "D48" "0.1680"
Set-TextValue "E48" "-24.00%"
Set-TextValue "E49" "0.12%"
Set-TextValue "E50" "0.12%"
